$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.992.56"
$ws.Range("E2").Value = "  -0.95%  "

# Row 3
$ws.Range("D3").Value = "1.760.40"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.87"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3890"
$ws.Range("E7").Value = "  +1.94%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3397"
$ws.Range("E8").Value = "  -1.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.32"
$ws.Range("E9").Value = "  -3.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  -2.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07206"
$ws.Range("E11").Value = "  -2.19%  "

# Row 12
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.28"
$ws.Range("E13").Value = "  -3.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.146"
$ws.Range("E14").Value = "  -4.65%  "

# Row 15
$ws.Range("D15").Value = "1.755.66"
$ws.Range("E15").Value = "  -1.91%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.049"
$ws.Range("E16").Value = "  -4.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("E17").Value = "  -1.69%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06602"
$ws.Range("E18").Value = "  -1.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.49"
$ws.Range("E19").Value = "  -2.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9977"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.91"
$ws.Range("E21").Value = "  -3.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  -3.89%  "

# Row 23
$ws.Range("D23").Value = "28.024.13"
$ws.Range("E23").Value = "  -0.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.62"
$ws.Range("E24").Value = "  -3.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.33"
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.93"
$ws.Range("E27").Value = "  -3.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.304"
$ws.Range("E28").Value = "  -4.61%  "

# Row 29
$ws.Range("D29").Value = "1.955.72"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.274"
$ws.Range("E30").Value = "  -12.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.27"
$ws.Range("E31").Value = "  -4.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  +3.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.815"
$ws.Range("E33").Value = "  -4.85%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08666"
$ws.Range("E34").Value = "  -2.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.04"
$ws.Range("E35").Value = "  -5.72%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06150"
$ws.Range("E36").Value = "  -3.25%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02277"
$ws.Range("E37").Value = "  -6.73%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.124"
$ws.Range("E38").Value = "  -3.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6470"
$ws.Range("E39").Value = "  -5.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2106"
$ws.Range("E40").Value = "  -3.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.498"
$ws.Range("E41").Value = "  +0.25%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.201"
$ws.Range("E42").Value = "  -3.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.850"
$ws.Range("E44").Value = "  -5.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("E45").Value = "  -3.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.825"
$ws.Range("E46").Value = "  -1.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("E47").Value = "  -5.11%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.45"
$ws.Range("E48").Value = "  -4.91%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.977"
$ws.Range("E49").Value = "  -5.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07003"
$ws.Range("E50").Value = "  -5.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.151"
$ws.Range("E51").Value = "  -4.49%  "
